$wb = $excel.ActiveWorkbook

# Row 32 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

# Row 33 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 173.73685
$ws.Range("I33").Value = 103.3125
$ws.Range("K33").Value = 103.3125
$ws.Range("M33").Value = 125.6875

# Row 40 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2057.1428
$ws.Range("I40").Value = 1850
$ws.Range("K40").Value = 1850
$ws.Range("M40").Value = -1675

# Row 132 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7582165.5
$ws.Range("I132").Value = 11911029
$ws.Range("J132").Value = 6654.3125
$ws.Range("K132").Value = 35733087
$ws.Range("L132").Value = 19962.9375
$ws.Range("M132").Value = -35730557
$ws.Range("N132").Value = -25022.9375

# Row 138 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1459.3334
$ws.Range("I138").Value = 647.7826
$ws.Range("J138").Value = 1704.9342
$ws.Range("K138").Value = 1943.3478
$ws.Range("L138").Value = 5114.8026
$ws.Range("M138").Value = 3196.6522
$ws.Range("N138").Value = -15394.8026

# Row 141 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 613.325
$ws.Range("I141").Value = 552.12823
$ws.Range("K141").Value = 1656.38469
$ws.Range("M141").Value = 3523.61531

# Row 43 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 7942.5
$ws.Range("J43").Value = 7942.5
$ws.Range("L43").Value = 7942.5
$ws.Range("N43").Value = -8568.5

# Row 109 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37774

# Row 134 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3750.0193
$ws.Range("I134").Value = 1252.5952
$ws.Range("J134").Value = 14239.2
$ws.Range("K134").Value = 3757.7856
$ws.Range("L134").Value = 42717.60000000001
$ws.Range("M134").Value = -1222.7856
$ws.Range("N134").Value = -47787.60000000001

# Row 58 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 894.4286
$ws.Range("I58").Value = 789.1818
$ws.Range("J58").Value = 1280.3334
$ws.Range("K58").Value = 789.1818
$ws.Range("L58").Value = 1280.3334
$ws.Range("M58").Value = -586.1818
$ws.Range("N58").Value = -1686.3334

# Row 132 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3065.132
$ws.Range("I132").Value = 3947.3057
$ws.Range("J132").Value = 1197
$ws.Range("K132").Value = 11841.9171
$ws.Range("L132").Value = 3591
$ws.Range("M132").Value = -9311.917099999999
$ws.Range("N132").Value = -8651

# Row 134 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1022.95746
$ws.Range("I134").Value = 1011.85364
$ws.Range("J134").Value = 1098.8334
$ws.Range("K134").Value = 3035.56092
$ws.Range("L134").Value = 3296.5002
$ws.Range("M134").Value = -500.5609199999999
$ws.Range("N134").Value = -8366.5002

# Row 136 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 894.4286
$ws.Range("I136").Value = 789.1818
$ws.Range("J136").Value = 1280.3334
$ws.Range("K136").Value = 2367.5454
$ws.Range("L136").Value = 3841.0002
$ws.Range("M136").Value = 182.4546
$ws.Range("N136").Value = -8941.0002

# Row 138 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 591780
$ws.Range("J138").Value = 591780
$ws.Range("L138").Value = 591780
$ws.Range("N138").Value = -602060

# Row 69 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3090.3333
$ws.Range("I69").Value = 999
$ws.Range("J69").Value = 3351.75
$ws.Range("K69").Value = 2997
$ws.Range("L69").Value = 10055.25
$ws.Range("M69").Value = -2186
$ws.Range("N69").Value = -11677.25

# Row 72 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 3090.3333
$ws.Range("I72").Value = 999
$ws.Range("J72").Value = 3351.75
$ws.Range("K72").Value = 8991
$ws.Range("L72").Value = 30165.75
$ws.Range("M72").Value = -4935
$ws.Range("N72").Value = -38277.75

# Row 113 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 729.25
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 729.25
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2187.75
$ws.Range("N113").Value = -6527.75
$ws.Range("M113").ClearContents()

# Row 122 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 834.5862
$ws.Range("J122").Value = 871.86365
$ws.Range("L122").Value = 7846.77285
$ws.Range("N122").Value = -12746.77285

# Row 131 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 28572856
$ws.Range("J131").Value = 1808
$ws.Range("L131").Value = 5424
$ws.Range("N131").Value = -15504

# Row 64 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 100 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 27166.666
$ws.Range("J100").Value = 27166.666
$ws.Range("L100").Value = 27166.666
$ws.Range("N100").Value = -29330.666

# Row 113 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1876.625
$ws.Range("I113").Value = 1880.4286
$ws.Range("J113").Value = 1850
$ws.Range("K113").Value = 1880.4286
$ws.Range("L113").Value = 1850
$ws.Range("M113").Value = 289.5714
$ws.Range("N113").Value = -6190

# Row 46 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6125
$ws.Range("J46").Value = 6125
$ws.Range("L46").Value = 6125
$ws.Range("N46").Value = -6501

# Row 122 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 12506579
$ws.Range("I122").Value = 15632593
$ws.Range("J122").Value = 2525
$ws.Range("K122").Value = 46897779
$ws.Range("L122").Value = 7575
$ws.Range("M122").Value = -46895329
$ws.Range("N122").Value = -12475

# Row 132 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 21119.576
$ws.Range("I132").Value = 1489.2593
$ws.Range("J132").Value = 42320.32
$ws.Range("K132").Value = 4467.7779
$ws.Range("L132").Value = 126960.96
$ws.Range("M132").Value = -1937.7779
$ws.Range("N132").Value = -132020.96

# Row 133 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 34099.855
$ws.Range("J133").Value = 34099.855
$ws.Range("L133").Value = 34099.855
$ws.Range("N133").Value = -39159.855

# Row 136 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4957.7407
$ws.Range("I136").Value = 5788.2856
$ws.Range("K136").Value = 17364.8568
$ws.Range("M136").Value = -14814.8568

# Row 132 on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2351.875
$ws.Range("I132").Value = 2245.3242
$ws.Range("J132").Value = 3666
$ws.Range("K132").Value = 6735.9726
$ws.Range("L132").Value = 10998
$ws.Range("M132").Value = -4205.9726
$ws.Range("N132").Value = -16058
